$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (Sponsor Work): add Friday (F13) hours, and daily total (I13)
$ws.Range("F13").Value = 1
$ws.Range("I13").Value = 1

# Row 14 (Daily Total / Weekly Total row): Friday column total, and updated weekly total
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 6

# Update the active selection to K13
$ws.Range("K13").Select()
